$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row index, Coin (B), Link (C), Price (D), Volume(1h) (E)
# A new row (OKB) was inserted at row 9, shifting the remaining coins down
# by one position; every row's Price/Volume(1h) reflects a refreshed quote.
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "27.509.32", "  -5.25%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.840.16", "  -4.40%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.000", "  -0.47%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "313.06", "  -3.79%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9995", "  -0.37%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4227", "  -7.87%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3631", "  -4.86%  "),
    @(9, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "44.05", "  -3.78%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07213", "  -6.93%  "),
    @(11, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.9006", "  -8.01%  "),
    @(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.59", "  -9.03%  "),
    @(13, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.853.31", "  -5.28%  "),
    @(14, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.579", "  -5.55%  "),
    @(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.295", "  -7.11%  "),
    @(16, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06794", "  -3.29%  "),
    @(17, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  -0.45%  "),
    @(18, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "77.30", "  -8.84%  "),
    @(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000008906", "  -6.25%  "),
    @(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9996", "  -0.40%  "),
    @(21, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "15.35", "  -8.10%  "),
    @(22, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "27.492.08", "  -5.38%  "),
    @(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.940", "  -7.72%  "),
    @(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "10.63", "  -3.70%  "),
    @(25, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.077.79", "  -3.24%  "),
    @(26, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.040", "  -1.19%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "151.42", "  -4.16%  "),
    @(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.22", "  -4.16%  "),
    @(29, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.288", "  -5.58%  "),
    @(30, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "110.74", "  -5.84%  "),
    @(31, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.692", "  -7.64%  "),
    @(32, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.08864", "  -4.94%  "),
    @(33, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7767", "  -9.76%  "),
    @(34, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.508", "  -11.53%  "),
    @(35, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.847", "  -5.53%  "),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.080", "  -13.18%  "),
    @(37, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9992", "  -0.44%  "),
    @(38, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05386", "  -5.34%  "),
    @(39, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.093", "  -5.01%  "),
    @(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.975", "  -4.20%  "),
    @(41, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01933", "  -5.44%  "),
    @(42, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.844", "  -7.98%  "),
    @(43, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.5037", "  -8.50%  "),
    @(44, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1636", "  -6.77%  "),
    @(45, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06613", "  -4.64%  "),
    @(46, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "8.211", "  -12.16%  "),
    @(47, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.4719", "  -8.90%  "),
    @(48, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "105.34", "  -4.78%  "),
    @(49, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "10.21", "  -9.06%  "),
    @(50, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9989", "  -0.40%  "),
    @(51, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.639", "  -6.90%  ")
)

foreach ($row in $rows) {
    $r = $row[0]

    # Leading apostrophe forces these (often numeric-looking) values to be
    # stored as literal text instead of being auto-converted to numbers;
    # resetting the style back to "Normal" afterwards drops the transient
    # quote-prefix formatting Excel applies so the cell keeps its original
    # (unstyled / General) appearance.
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
    $ws.Cells.Item($r, 5).Style = "Normal"
}
